$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7517.0713
$ws.Range("J28").Value = 7272.778
$ws.Range("L28").Value = 7272.778
$ws.Range("N28").Value = -8242.778
$ws.Range("H40").Value = 2397.5
$ws.Range("J40").Value = 2400
$ws.Range("L40").Value = 2400
$ws.Range("N40").Value = -2750
$ws.Range("H62").Value = 8430.571
$ws.Range("I62").Value = 8499.5
$ws.Range("K62").Value = 8499.5
$ws.Range("M62").Value = -7875.5
$ws.Range("H65").Value = 8430.571
$ws.Range("I65").Value = 8499.5
$ws.Range("K65").Value = 42497.5
$ws.Range("M65").Value = -39377.5
$ws.Range("H137").Value = 1754.25
$ws.Range("I137").Value = 1454.5428
$ws.Range("K137").Value = 4363.6284
$ws.Range("M137").Value = -1813.6284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 915.2414
$ws.Range("I2").Value = 881.12964
$ws.Range("J2").Value = 1375.75
$ws.Range("K2").Value = 881.12964
$ws.Range("L2").Value = 1375.75
$ws.Range("M2").Value = -768.12964
$ws.Range("N2").Value = -1601.75
$ws.Range("H45").Value = 2971.3333
$ws.Range("I45").Value = 2969.2856
$ws.Range("K45").Value = 2969.2856
$ws.Range("M45").Value = -2592.2856
$ws.Range("H61").Value = 3581.5
$ws.Range("I61").Value = 2924.8
$ws.Range("K61").Value = 2924.8
$ws.Range("M61").Value = -2712.8
$ws.Range("H102").Value = 23491.637
$ws.Range("I102").Value = 15840.8
$ws.Range("K102").Value = 15840.8
$ws.Range("M102").Value = -14218.8
$ws.Range("H116").Value = 915.2414
$ws.Range("I116").Value = 881.12964
$ws.Range("J116").Value = 1375.75
$ws.Range("K116").Value = 881.12964
$ws.Range("L116").Value = 1375.75
$ws.Range("M116").Value = 1412.87036
$ws.Range("N116").Value = -5963.75
$ws.Range("H122").Value = 11658.333
$ws.Range("I122").Value = 12520.667
$ws.Range("K122").Value = 37562.001
$ws.Range("M122").Value = -35112.001
$ws.Range("H132").Value = 2540.5186
$ws.Range("I132").Value = 2263.96
$ws.Range("J132").Value = 5997.5
$ws.Range("K132").Value = 6791.88
$ws.Range("L132").Value = 17992.5
$ws.Range("M132").Value = -4261.88
$ws.Range("N132").Value = -23052.5
$ws.Range("H136").Value = 3581.5
$ws.Range("I136").Value = 2924.8
$ws.Range("K136").Value = 8774.400000000001
$ws.Range("M136").Value = -6224.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 915.2414
$ws.Range("I3").Value = 881.12964
$ws.Range("J3").Value = 1375.75
$ws.Range("K3").Value = 881.12964
$ws.Range("L3").Value = 1375.75
$ws.Range("M3").Value = -767.12964
$ws.Range("N3").Value = -1603.75
$ws.Range("H86").Value = 13891264
$ws.Range("I86").Value = 13891264
$ws.Range("K86").Value = 13891264
$ws.Range("M86").Value = -13890141
$ws.Range("H89").Value = 13891264
$ws.Range("I89").Value = 13891264
$ws.Range("K89").Value = 69456320
$ws.Range("M89").Value = -69450704
$ws.Range("H99").Value = 4109.2
$ws.Range("I99").Value = 4188.4287
$ws.Range("K99").Value = 4188.4287
$ws.Range("M99").Value = -2690.4287
$ws.Range("H105").Value = 3866.889
$ws.Range("I105").Value = 3601.25
$ws.Range("K105").Value = 3601.25
$ws.Range("M105").Value = -1854.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22232.215
$ws.Range("I16").Value = 16275.2
$ws.Range("J16").Value = 37124.75
$ws.Range("K16").Value = 16275.2
$ws.Range("L16").Value = 37124.75
$ws.Range("M16").Value = -15988.2
$ws.Range("N16").Value = -37698.75
$ws.Range("H94").Value = 616.75
$ws.Range("I94").Value = 1029.75
$ws.Range("J94").Value = 410.25
$ws.Range("K94").Value = 1029.75
$ws.Range("L94").Value = 410.25
$ws.Range("M94").Value = -578.75
$ws.Range("N94").Value = -1312.25
$ws.Range("H113").Value = 22232.215
$ws.Range("I113").Value = 16275.2
$ws.Range("J113").Value = 37124.75
$ws.Range("K113").Value = 16275.2
$ws.Range("L113").Value = 37124.75
$ws.Range("M113").Value = -14105.2
$ws.Range("N113").Value = -41464.75
$ws.Range("H122").Value = 115974.61
$ws.Range("I122").Value = 144096.1
$ws.Range("J122").Value = 12862.5
$ws.Range("K122").Value = 432288.3
$ws.Range("L122").Value = 38587.5
$ws.Range("M122").Value = -429838.3
$ws.Range("N122").Value = -43487.5
$ws.Range("H132").Value = 2280.8262
$ws.Range("I132").Value = 1973.6
$ws.Range("K132").Value = 5920.799999999999
$ws.Range("M132").Value = -3390.799999999999
$ws.Range("H134").Value = 1064.1666
$ws.Range("I134").Value = 1046.6666
$ws.Range("J134").Value = 1099.1666
$ws.Range("K134").Value = 3139.9998
$ws.Range("L134").Value = 3297.4998
$ws.Range("M134").Value = -604.9998000000001
$ws.Range("N134").Value = -8367.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 215.6
$ws.Range("I12").Value = 106.28571
$ws.Range("J12").Value = 311.25
$ws.Range("K12").Value = 318.85713
$ws.Range("L12").Value = 933.75
$ws.Range("M12").Value = -145.85713
$ws.Range("N12").Value = -1279.75
$ws.Range("H92").Value = 774.3333
$ws.Range("J92").Value = 774.3333
$ws.Range("L92").Value = 2322.9999
$ws.Range("N92").Value = -4818.9999
$ws.Range("H136").Value = 6423.2046
$ws.Range("I136").Value = 4816.864
$ws.Range("K136").Value = 14450.592
$ws.Range("M136").Value = -9350.591999999999
$ws.Range("H137").Value = 4329.6313
$ws.Range("J137").Value = 3736.9092
$ws.Range("L137").Value = 11210.7276
$ws.Range("N137").Value = -21410.7276
$ws.Range("H140").Value = 1448.5
$ws.Range("I140").Value = 970.85
$ws.Range("K140").Value = 2912.55
$ws.Range("M140").Value = 2267.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6629.73
$ws.Range("I132").Value = 6024.46
$ws.Range("K132").Value = 18073.38
$ws.Range("M132").Value = -15543.38

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7249.5
$ws.Range("I7").Value = 7000
$ws.Range("K7").Value = 7000
$ws.Range("M7").Value = -6888
$ws.Range("H80").Value = 79500
$ws.Range("I80").Value = 79000
$ws.Range("K80").Value = 79000
$ws.Range("M80").Value = -77877
$ws.Range("H83").Value = 79500
$ws.Range("I83").Value = 79000
$ws.Range("K83").Value = 237000
$ws.Range("M83").Value = -231384
$ws.Range("H126").Value = 7249.5
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -18530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4715.1304
$ws.Range("I126").Value = 4617.381
$ws.Range("K126").Value = 13852.143
$ws.Range("M126").Value = -11382.143
$ws.Range("H132").Value = 6785.964
$ws.Range("I132").Value = 5563.3784
$ws.Range("J132").Value = 9299.056
$ws.Range("K132").Value = 16690.1352
$ws.Range("L132").Value = 27897.168
$ws.Range("M132").Value = -14160.1352
$ws.Range("N132").Value = -32957.16800000001
